$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# The task table has a header row (row 1) followed by rows for tasks
# 1..13. Several "ASSIGN" cells (last column) were left blank and now
# get the responsible person's name; the "DETAILS" cell for task 8 had
# its sentence split across five runs that get collapsed into one.

# Task 4 -> ASSIGN: Hải
$t.Cell(5, 4).Range.Text = "Hải"

# Task 5 -> ASSIGN: Vinh
$t.Cell(6, 4).Range.Text = "Vinh"

# Task 6 -> ASSIGN: Hải - Vinh
$t.Cell(7, 4).Range.Text = "Hải - Vinh"

# Task 7 -> ASSIGN: Phúc
$t.Cell(8, 4).Range.Text = "Phúc"

# Task 8 -> DETAILS: merge the five runs into a single run with the
# same overall sentence (no leading/trailing space differences).
$detailsCell = $t.Cell(9, 2)
$detailsCell.Range.Find.Execute(
    "Xây dựng module cho Digi ME 9210 để truyền nhận data thông qua UART.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Xây dựng module cho Digi ME 9210 để truyền nhận data thông qua UART.",
    2)

# Task 8 -> ASSIGN: Phúc
$t.Cell(9, 4).Range.Text = "Phúc"

# Task 9 -> ASSIGN: Team
$t.Cell(10, 4).Range.Text = "Team"

# Task 10 -> ASSIGN: Vinh
$t.Cell(11, 4).Range.Text = "Vinh"

Write-Host "Task list ASSIGN column and deadline details updated."
